$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# This edit inserts a new "latest price" row (Sl.no. 21, 02-11-2025)
# right above the existing "16 / 01-11-2025" row (current row 6),
# pushing every row below it down by one. Row 1 (header) and rows
# 2-5 are untouched.
# ------------------------------------------------------------------

# Hyperlinks in this engine don't automatically re-target when rows
# shift, so drop them all up-front and rebuild the full F2:F22 set
# afterwards in the correct, final row order.
$ws.Hyperlinks.Delete()

# Shift existing data rows 6:21 down to 7:22, carrying values/styles.
$ws.Rows("6:6").Insert()

# Make sure the new date cell is stored as literal text (not an
# auto-converted Excel date serial) before assigning it.
$ws.Range("E6").NumberFormat = "@"

# New row 6 content.
$ws.Range("A6").Value = 21
$ws.Range("B6").Value = "ALUMINIUM INGOT"
$ws.Range("C6").Value = "IE07"
$ws.Range("D6").Value = 296.05
$ws.Range("E6").Value = "02-11-2025"
$ws.Range("F6").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Re-apply the standard data-row formatting (the inserted row above
# picked it up already; this just normalizes the NumberFormat="@"
# tweak back to the shared style used by every other data row).
$ws.Range("A7:F7").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)

# Rebuild every hyperlink, F2 through F22, pointing at the URL text
# already sitting in each cell.
$ws.Hyperlinks.Add($ws.Range("F2"),  $ws.Range("F2").Value2)
$ws.Hyperlinks.Add($ws.Range("F3"),  $ws.Range("F3").Value2)
$ws.Hyperlinks.Add($ws.Range("F4"),  $ws.Range("F4").Value2)
$ws.Hyperlinks.Add($ws.Range("F5"),  $ws.Range("F5").Value2)
$ws.Hyperlinks.Add($ws.Range("F6"),  $ws.Range("F6").Value2)
$ws.Hyperlinks.Add($ws.Range("F7"),  $ws.Range("F7").Value2)
$ws.Hyperlinks.Add($ws.Range("F8"),  $ws.Range("F8").Value2)
$ws.Hyperlinks.Add($ws.Range("F9"),  $ws.Range("F9").Value2)
$ws.Hyperlinks.Add($ws.Range("F10"), $ws.Range("F10").Value2)
$ws.Hyperlinks.Add($ws.Range("F11"), $ws.Range("F11").Value2)
$ws.Hyperlinks.Add($ws.Range("F12"), $ws.Range("F12").Value2)
$ws.Hyperlinks.Add($ws.Range("F13"), $ws.Range("F13").Value2)
$ws.Hyperlinks.Add($ws.Range("F14"), $ws.Range("F14").Value2)
$ws.Hyperlinks.Add($ws.Range("F15"), $ws.Range("F15").Value2)
$ws.Hyperlinks.Add($ws.Range("F16"), $ws.Range("F16").Value2)
$ws.Hyperlinks.Add($ws.Range("F17"), $ws.Range("F17").Value2)
$ws.Hyperlinks.Add($ws.Range("F18"), $ws.Range("F18").Value2)
$ws.Hyperlinks.Add($ws.Range("F19"), $ws.Range("F19").Value2)
$ws.Hyperlinks.Add($ws.Range("F20"), $ws.Range("F20").Value2)
$ws.Hyperlinks.Add($ws.Range("F21"), $ws.Range("F21").Value2)
$ws.Hyperlinks.Add($ws.Range("F22"), $ws.Range("F22").Value2)

# Adding hyperlinks applies Excel's built-in "Hyperlink" style
# (blue/underline) to each cell; restore the plain data-row style
# used throughout the sheet without touching the stored text/links.
$ws.Range("E2").Copy()
$ws.Range("F2:F22").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false
